$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This CI "handback status" report gets a new row appended to each of its
# three tables (Overview, zh-cn, de-de) for a freshly-generated handback
# file (82304c21-...) plus its dependent file (d8a125cf-...). The existing
# row 2 on the zh-cn / de-de sheets is rewritten in place to reflect the
# newest handoff/handback run (new guid, new xliff hash, refreshed
# timestamps) while sheet1 (Overview) keeps row2 untouched and only gains a
# new row3 for the new markdown file.
# ---------------------------------------------------------------------------

function Set-Text($ws, $a1, $text) {
    # Force text interpretation even for strings that look like booleans /
    # numbers / dates, by prefixing with an apostrophe (Excel's "treat as
    # text" marker) - this never ends up stored in the cell itself.
    $ws.Range($a1).Value2 = "'" + $text
}

# =====================  Sheet 1: "Overview"  ================================
$ws1 = $wb.Worksheets.Item("Overview")

$lo1 = $ws1.ListObjects.Item(1)
$lo1.ListRows.Add() | Out-Null

Set-Text $ws1 "A3" "d8a125cf-0be1-4789-a552-c4854026fd03.md"
Set-Text $ws1 "B3" "e2e\d8a125cf-0be1-4789-a552-c4854026fd03.md"
Set-Text $ws1 "C3" ".md"
Set-Text $ws1 "E3" "Handed back: in sync with en-US"
Set-Text $ws1 "F3" "Handed back: in sync with en-US"
Set-Text $ws1 "G3" "2016-09-07 03:17:56"
$ws1.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"


# Refresh row2's existing markdown-file hyperlink text/target (the handback
# filename moved from 4db14c90-... to 82304c21-...). This engine's COM
# surface doesn't let us mutate a Hyperlink's Address/TextToDisplay in
# place, so drop the old one and re-add - deleting first frees up its rId
# so the replacement reclaims the same id.
Set-Text $ws1 "A2" "82304c21-bfe5-443d-8050-97c14a7bf7ad.md"
Set-Text $ws1 "B2" "e2e\82304c21-bfe5-443d-8050-97c14a7bf7ad.md"
Set-Text $ws1 "G2" "2016-09-07 03:17:56"
$ws1.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws1.Range("B2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1c1397601d9442fc55afa3b5232f8d0ef3afdba/e2e/82304c21-bfe5-443d-8050-97c14a7bf7ad.md", "", "", "e2e\82304c21-bfe5-443d-8050-97c14a7bf7ad.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1c1397601d9442fc55afa3b5232f8d0ef3afdba/e2e/d8a125cf-0be1-4789-a552-c4854026fd03.md", "", "", "e2e\d8a125cf-0be1-4789-a552-c4854026fd03.md") | Out-Null

# =====================  Sheet 2: "zh-cn"  ====================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$lo2 = $ws2.ListObjects.Item(1)

# --- row 2 (latest run for the primary handback file) ----------------------
Set-Text $ws2 "A2" "82304c21-bfe5-443d-8050-97c14a7bf7ad.md"
Set-Text $ws2 "B2" ".md"
Set-Text $ws2 "C2" "Handed back: in sync with en-US"
Set-Text $ws2 "D2" "e2e"
Set-Text $ws2 "E2" "ht"
Set-Text $ws2 "F2" "False"
Set-Text $ws2 "G2" "82304c21-bfe5-443d-8050-97c14a7bf7ad.fbe15d2db5440491d03be5bfb17e966711c77855.zh-cn.xlf"
Set-Text $ws2 "H2" "2016-09-07 03:17:51"
$ws2.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws2 "I2" "82304c21-bfe5-443d-8050-97c14a7bf7ad.md"
Set-Text $ws2 "J2" "82304c21-bfe5-443d-8050-97c14a7bf7ad.fbe15d2db5440491d03be5bfb17e966711c77855.zh-cn.xlf"
Set-Text $ws2 "K2" "2016-09-07 03:18:19"
$ws2.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws2 "L2" ""
Set-Text $ws2 "M2" "True"
Set-Text $ws2 "N2" ""
Set-Text $ws2 "O2" "False"
Set-Text $ws2 "P2" ""

# Drop + re-add (in the same ref order) so the freed rId2/rId3 get reclaimed
# by the refreshed links instead of appending new ones.
$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Range("I2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1c1397601d9442fc55afa3b5232f8d0ef3afdba/e2e/82304c21-bfe5-443d-8050-97c14a7bf7ad.md", "", "", "82304c21-bfe5-443d-8050-97c14a7bf7ad.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4729df56f0222146f2442925d8cd8396e0c24240/e2e/82304c21-bfe5-443d-8050-97c14a7bf7ad.md", "", "", "82304c21-bfe5-443d-8050-97c14a7bf7ad.md") | Out-Null

# --- row 3 (new dependent file) ---------------------------------------------
$lo2.ListRows.Add() | Out-Null

Set-Text $ws2 "A3" "d8a125cf-0be1-4789-a552-c4854026fd03.md"
Set-Text $ws2 "B3" ".md"
Set-Text $ws2 "C3" "Handed back: in sync with en-US"
Set-Text $ws2 "D3" "e2e"
Set-Text $ws2 "E3" "ht"
Set-Text $ws2 "F3" "True"
Set-Text $ws2 "G3" "d8a125cf-0be1-4789-a552-c4854026fd03.c81b41df33f2b93fc72300653e0bb6c6eddc303a.zh-cn.xlf"
Set-Text $ws2 "H3" "2016-09-07 03:17:51"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws2 "I3" "d8a125cf-0be1-4789-a552-c4854026fd03.md"
Set-Text $ws2 "J3" "d8a125cf-0be1-4789-a552-c4854026fd03.c81b41df33f2b93fc72300653e0bb6c6eddc303a.zh-cn.xlf"
Set-Text $ws2 "K3" "2016-09-07 03:18:19"
$ws2.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws2 "L3" ""
Set-Text $ws2 "M3" "True"
Set-Text $ws2 "N3" ""
Set-Text $ws2 "O3" "False"
Set-Text $ws2 "P3" ""

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1c1397601d9442fc55afa3b5232f8d0ef3afdba/e2e/d8a125cf-0be1-4789-a552-c4854026fd03.md", "", "", "d8a125cf-0be1-4789-a552-c4854026fd03.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4729df56f0222146f2442925d8cd8396e0c24240/e2e/d8a125cf-0be1-4789-a552-c4854026fd03.md", "", "", "d8a125cf-0be1-4789-a552-c4854026fd03.md") | Out-Null

# =====================  Sheet 3: "de-de"  ====================================
$ws3 = $wb.Worksheets.Item("de-de")

$lo3 = $ws3.ListObjects.Item(1)

# --- row 2 (latest run for the primary handback file) ----------------------
Set-Text $ws3 "A2" "82304c21-bfe5-443d-8050-97c14a7bf7ad.md"
Set-Text $ws3 "B2" ".md"
Set-Text $ws3 "C2" "Handed back: in sync with en-US"
Set-Text $ws3 "D2" "e2e"
Set-Text $ws3 "E2" "ht"
Set-Text $ws3 "F2" "False"
Set-Text $ws3 "G2" "82304c21-bfe5-443d-8050-97c14a7bf7ad.fbe15d2db5440491d03be5bfb17e966711c77855.de-de.xlf"
Set-Text $ws3 "H2" "2016-09-07 03:17:56"
$ws3.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws3 "I2" "82304c21-bfe5-443d-8050-97c14a7bf7ad.md"
Set-Text $ws3 "J2" "82304c21-bfe5-443d-8050-97c14a7bf7ad.fbe15d2db5440491d03be5bfb17e966711c77855.de-de.xlf"
Set-Text $ws3 "K2" "2016-09-07 03:18:27"
$ws3.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws3 "L2" ""
Set-Text $ws3 "M2" "True"
Set-Text $ws3 "N2" ""
Set-Text $ws3 "O2" "False"
Set-Text $ws3 "P2" ""

# Drop + re-add (in the same ref order) so the freed rId2/rId3 get reclaimed
# by the refreshed links instead of appending new ones.
$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Range("I2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1c1397601d9442fc55afa3b5232f8d0ef3afdba/e2e/82304c21-bfe5-443d-8050-97c14a7bf7ad.md", "", "", "82304c21-bfe5-443d-8050-97c14a7bf7ad.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c5c745f00750db7b0b558c0b6d6285c42a390ab0/e2e/82304c21-bfe5-443d-8050-97c14a7bf7ad.md", "", "", "82304c21-bfe5-443d-8050-97c14a7bf7ad.md") | Out-Null

# --- row 3 (new dependent file) ---------------------------------------------
$lo3.ListRows.Add() | Out-Null

Set-Text $ws3 "A3" "d8a125cf-0be1-4789-a552-c4854026fd03.md"
Set-Text $ws3 "B3" ".md"
Set-Text $ws3 "C3" "Handed back: in sync with en-US"
Set-Text $ws3 "D3" "e2e"
Set-Text $ws3 "E3" "ht"
Set-Text $ws3 "F3" "True"
Set-Text $ws3 "G3" "d8a125cf-0be1-4789-a552-c4854026fd03.c81b41df33f2b93fc72300653e0bb6c6eddc303a.de-de.xlf"
Set-Text $ws3 "H3" "2016-09-07 03:17:56"
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws3 "I3" "d8a125cf-0be1-4789-a552-c4854026fd03.md"
Set-Text $ws3 "J3" "d8a125cf-0be1-4789-a552-c4854026fd03.c81b41df33f2b93fc72300653e0bb6c6eddc303a.de-de.xlf"
Set-Text $ws3 "K3" "2016-09-07 03:18:27"
$ws3.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-Text $ws3 "L3" ""
Set-Text $ws3 "M3" "True"
Set-Text $ws3 "N3" ""
Set-Text $ws3 "O3" "False"
Set-Text $ws3 "P3" ""

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1c1397601d9442fc55afa3b5232f8d0ef3afdba/e2e/d8a125cf-0be1-4789-a552-c4854026fd03.md", "", "", "d8a125cf-0be1-4789-a552-c4854026fd03.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c5c745f00750db7b0b558c0b6d6285c42a390ab0/e2e/d8a125cf-0be1-4789-a552-c4854026fd03.md", "", "", "d8a125cf-0be1-4789-a552-c4854026fd03.md") | Out-Null

Write-Host "Generate Report for Handback: done"
